# Remove the "platform" field from the ENA experiment metadata template.
# "platform" is non-mandatory once "instrument_model" is supplied, so:
#   - drop the "platform" column (and its values) from the "experiment" sheet
#   - drop the "platform" column (and its values) from the "cv_experiment" lookup sheet
#   - remove the "platform" named range
#   - point the "instrumentmodel" named range at its new column

$wb = $excel.ActiveWorkbook

# 1. experiment sheet: column M is "platform" (column N, "instrument_model",
#    shifts left into M and picks up its data validation automatically).
$wsExperiment = $wb.Worksheets.Item("experiment")
$wsExperiment.Columns("M").Delete()

# 2. cv_experiment sheet: column M held the platform value list
#    (LS454, ILLUMINA, ...); column N held the instrument_model list and
#    shifts left into M.
$wsCv = $wb.Worksheets.Item("cv_experiment")
$wsCv.Columns("M").Delete()

# 3. Defined names: drop "platform" entirely, and repoint "instrumentmodel"
#    at its new column M (was N).
$wb.Names.Item("platform").Delete()
$wb.Names.Item("instrumentmodel").RefersTo = "='cv_experiment'!`$M`$1:`$M`$83"
